# "Generate Report for Handoff"
# The localization status has moved from "Handed back: in sync with en-US"
# to "Ready for handoff", and the associated timestamps were refreshed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-04 23:04:37"
$overview.Columns("E").ColumnWidth = 16.38265482584637
$overview.Columns("F").ColumnWidth = 16.38265482584637

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-04 23:04:32"
$zhcn.Columns("C").ColumnWidth = 16.38265482584637

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-04 23:04:37"
$dede.Columns("C").ColumnWidth = 16.38265482584637
